$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert new "lang_code" column, shift rest left, rename ---
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# --- Data rows ---
$data = @(
    @("eng", "FNR", "Fingerprint",           "Finger prints of the applicant",     $true),
    @("eng", "IRS", "Iris",                  "Iris of the applicant",              $true),
    @("eng", "PHT", "Photo",                 "Photo of the face of the applicant", $true),
    @("fra", "FNR", "Empreintes digitales",  "Empreintes digitales du demandeur",  $true),
    @("fra", "IRS", "Iris",                  "Iris du demandeur",                  $true),
    @("fra", "PHT", "Photo",                 "Photo du visage du demandeur",       $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# --- Apply the header style (index s="1") to column A of each data row,
#     matching the workbook's original "first column" emphasis style ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
